$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 108, shifting existing rows 108:126 down to 109:127.
$ws.Rows(108).Insert()

# Populate the new row 108 with the latest weekly price record.
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = 'Vega Modelo de Temuco'
$ws.Range("C108").Value = 'La Araucanía'
$ws.Range("D108").Value = 44637
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100112012
$ws.Range("G108").Value = 'Espinaca'
$ws.Range("H108").Value = 'Sin especificar'
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 50
$ws.Range("K108").Value = 10000
$ws.Range("L108").Value = 10000
$ws.Range("M108").Value = 10000
$ws.Range("N108").Value = '$/docena de atados'
$ws.Range("O108").Value = 'Región de La Araucanía'
$ws.Range("P108").Value = 3333
$ws.Range("Q108").Value = 3
$ws.Range("R108").Value = 'Hortaliza'
